{"js": "// Adds the new \"UML\" section (heading + bulleted list) right after the\n// paragraph \"Las subclases heredan los miembros de la superclase - metodos\n// y atributos.\" and before the trailing empty list paragraph at the end\n// of the document.\n//\n// Strategy: every new paragraph is inserted with InsertLocation.after\n// directly on the anchor paragraph (never chained after the\n// previously-inserted one). That way each insertion always inherits the\n// anchor's own formatting (ListParagraph / bullet list numId=1, ilvl=0),\n// which is exactly what the two top-level bullets need \"for free\", and is\n// the correct starting point for the others (the heading then has its list\n// formatting removed, and the sub-bullets get promoted to level 1).\n// Building the run in reverse order - last sentence first - means\n// everything still reads top-to-bottom afterwards.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet anchor = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Las subclases heredan\") === 0) {\n    anchor = items[i];\n  }\n}\nif (!anchor) {\n  anchor = items[items.length - 2];\n}\n\n// 6) Implementaci\u00f3n: se codifica en un lenguaje OO.  (sub-bullet, level 1)\nconst p6 = anchor.insertParagraph(\n  \"Implementaci\u00f3n: se codifica en un lenguaje OO.\",\n  Word.InsertLocation.after\n);\np6.listItemOrNullObject.level = 1;\n\n// 5) Dise\u00f1o: Como es que los objetos se relacionan pata dar una soluci\u00f3n al problema.  (sub-bullet, level 1)\nconst p5 = anchor.insertParagraph(\n  \"Dise\u00f1o: Como es que los objetos se relacionan pata dar una soluci\u00f3n al problema.\",\n  Word.InsertLocation.after\n);\np5.listItemOrNullObject.level = 1;\n\n// 4) An\u00e1lisis: identifica objetos, as\u00ed como sus atributos y actividades u operaciones.  (sub-bullet, level 1)\nconst p4 = anchor.insertParagraph(\n  \"An\u00e1lisis: identifica objetos, as\u00ed como sus atributos y actividades u operaciones.\",\n  Word.InsertLocation.after\n);\np4.listItemOrNullObject.level = 1;\n\n// 3) Es una herramienta para el an\u00e1lisis, dise\u00f1o y la implementaci\u00f3n de sistemas.  (top-level bullet, level 0)\nconst p3 = anchor.insertParagraph(\n  \"Es una herramienta para el an\u00e1lisis, dise\u00f1o y la implementaci\u00f3n de sistemas.\",\n  Word.InsertLocation.after\n);\n\n// 2) El lenguaje unificado de modelado es un est\u00e1ndar para el an\u00e1lisis y desarrollo orientado a objetos.  (top-level bullet, level 0)\nconst p2 = anchor.insertParagraph(\n  \"El lenguaje unificado de modelado es un est\u00e1ndar para el an\u00e1lisis y desarrollo orientado a objetos.\",\n  Word.InsertLocation.after\n);\n\n// 1) UML  (plain paragraph heading - no list formatting)\nconst p1 = anchor.insertParagraph(\"\", Word.InsertLocation.after);\np1.detachFromList();\np1.style = \"Normal\";\nawait context.sync();\n\np1.insertText(\"UML\", Word.InsertLocation.end);\np1.alignment = Word.Alignment.justified;\np1.lineSpacing = 13.8; // 276 twips / 20 = 13.8 pt\n\nawait context.sync();\n", "ps1": "# Adds the new \"UML\" section (heading + bulleted list) right after the\n# paragraph \"Las subclases heredan los miembros de la superclase - metodos\n# y atributos.\" and before the trailing empty list paragraph at the end\n# of the document.\n#\n# Strategy: every new paragraph is inserted immediately after the anchor\n# paragraph (not chained after the previously-inserted one). That way each\n# insertion always inherits the anchor's own formatting (ListParagraph /\n# bullet list numId=1, ilvl=0), which is exactly what the two top-level\n# bullets need, and is the correct starting point for the others (which we\n# then touch up: the heading loses the list formatting, the sub-bullets get\n# promoted to ilvl=1). Building the run in reverse order - last sentence\n# first - means everything still reads top-to-bottom afterwards.\n\n$d = $word.ActiveDocument\n\n# Find the \"Las subclases...\" paragraph by its text so this keeps working\n# even if paragraph indices shift.\n$anchor = $null\nforeach ($para in $d.Paragraphs) {\n    if ($para.Range.Text -like \"Las subclases heredan*\") {\n        $anchor = $para\n    }\n}\nif ($anchor -eq $null) {\n    $anchor = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n}\n\n# 6) Implementaci\u00f3n: se codifica en un lenguaje OO.  (sub-bullet, ilvl=1)\n$anchor.Range.InsertParagraphAfter()\n$n6 = $anchor.Next()\n$n6.Range.InsertAfter(\"Implementaci\u00f3n: se codifica en un lenguaje OO.\")\n$n6.Range.ListFormat.ListLevelNumber = 2\n\n# 5) Dise\u00f1o: Como es que los objetos se relacionan pata dar una soluci\u00f3n al problema.  (sub-bullet, ilvl=1)\n$anchor.Range.InsertParagraphAfter()\n$n5 = $anchor.Next()\n$n5.Range.InsertAfter(\"Dise\u00f1o: Como es que los objetos se relacionan pata dar una soluci\u00f3n al problema.\")\n$n5.Range.ListFormat.ListLevelNumber = 2\n\n# 4) An\u00e1lisis: identifica objetos, as\u00ed como sus atributos y actividades u operaciones.  (sub-bullet, ilvl=1)\n$anchor.Range.InsertParagraphAfter()\n$n4 = $anchor.Next()\n$n4.Range.InsertAfter(\"An\u00e1lisis: identifica objetos, as\u00ed como sus atributos y actividades u operaciones.\")\n$n4.Range.ListFormat.ListLevelNumber = 2\n\n# 3) Es una herramienta para el an\u00e1lisis, dise\u00f1o y la implementaci\u00f3n de sistemas.  (top-level bullet, ilvl=0)\n$anchor.Range.InsertParagraphAfter()\n$n3 = $anchor.Next()\n$n3.Range.InsertAfter(\"Es una herramienta para el an\u00e1lisis, dise\u00f1o y la implementaci\u00f3n de sistemas.\")\n\n# 2) El lenguaje unificado de modelado es un est\u00e1ndar para el an\u00e1lisis y desarrollo orientado a objetos.  (top-level bullet, ilvl=0)\n$anchor.Range.InsertParagraphAfter()\n$n2 = $anchor.Next()\n$n2.Range.InsertAfter(\"El lenguaje unificado de modelado es un est\u00e1ndar para el an\u00e1lisis y desarrollo orientado a objetos.\")\n\n# 1) UML  (plain paragraph heading - no list formatting)\n$anchor.Range.InsertParagraphAfter()\n$n1 = $anchor.Next()\n$n1.Range.ListFormat.RemoveNumbers()\n$n1.Style = \"Normal\"\n$n1.Range.InsertAfter(\"UML\")\n$n1.Alignment = 3            # wdAlignParagraphJustify\n$n1.LineSpacingRule = 5      # wdLineSpaceMultiple\n$n1.LineSpacing = 13.8       # 276 twips / 20 = 13.8 pt\n\nWrite-Output \"done\"\n"}
